$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline approach: for numeric-looking text in the Price column (D),
# force text storage (matches source data which keeps these as text/inlineStr)
# by briefly applying a Text number format, then clearing the style so no stray
# formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.230.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.90%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.575.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.49%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "517.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.33%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.588.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.68%  "

$ws.Range("E11").Value = "  -4.54%  "

$ws.Range("E12").Value = "  -2.55%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.133"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.030.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.176.69"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.563.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.29%  "

$ws.Range("E18").Value = "  -3.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "335.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.95%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.45%  "

$ws.Range("E25").Value = "  -1.62%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("B27").Value = "Polygon"
$ws.Range("C27").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.400"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.08%  "

$ws.Range("E29").Value = "  -0.04%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0707"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -11.19%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.95%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.86%  "

$ws.Range("E33").Value = "  -3.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.62%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.89"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.47%  "

$ws.Range("E36").Value = "  -5.46%  "

$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.829"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("E40").Value = "  -2.62%  "

$ws.Range("E41").Value = "  -3.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "271.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.588"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0941"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0513"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.72%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.970.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.18%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0218"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.08%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.72%  "
